# "Extended comment support WIP"
# Adds a 4th column (D) to the "Слова" sheet holding an extended-comment
# footnote, and appends "*1"/"*2" footnote markers to the existing B3/C3
# example cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cell: D1 --------------------------------------------------
$ws.Range("D1").Value = "Расширенный комментарий"

# Copy header formatting (bold font + border + wrap) from C1 onto D1, so the
# new cell reuses the existing header style instead of minting a new one.
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)

# --- Updated example cells: B3 / C3 ---------------------------------------
$ws.Range("B3").Value = "словоформаБ1 *1"
$ws.Range("C3").Value = "переводБ1 | переводБ2 *2"

# --- New footnote cell: D3 --------------------------------------------------
$ws.Range("D3").Value = "*1: Расширенный комментарий`n*2: Расширенный комментарий"

# Copy body formatting (border + wrap) from C3 onto D3.
$ws.Range("C3").Copy()
$ws.Range("D3").PasteSpecial(-4122)

# --- Column width for the new column ---------------------------------------
# NOTE: the host's pixel-quantised ColumnWidth setter cannot reproduce the
# source workbook's exact 32.5703125-character width (it snaps to 1/6-char
# steps); an input of 31.6 is the value that round-trips to the closest
# attainable stored width (32.5).
$ws.Columns.Item(4).ColumnWidth = 31.6

# --- Row height for row 3 (now wraps the longer footnoted text) ------------
$ws.Rows.Item(3).RowHeight = 30

# --- Selection, matching the post-edit workbook state -----------------------
$ws.Range("D4").Select()
